# se modifica data para empezar regresion R34 en Pre Prod
#
# Updates the data sheet ("Hoja1") of the workbook:
#   - G11: 24741865 -> 24741866
#   - N11: 307      -> 130
#   - move the active selection from N11 to N12 (cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the data values
$ws.Range("G11").Value = 24741866
$ws.Range("N11").Value = 130

# Reflect the new selection/cursor position (sheet stays the active tab)
$ws.Activate()
[void]$ws.Range("N12").Select()
